# The document's headers/footers each carry one inline picture:
#   - the Pearson logo (alt text ends in "PearsonLogo.png"), whose
#     wp:docPr/pic:cNvPr "name" attribute should read "image2.png"
#     (it currently reads "image1.png")
#   - the BTec logo (alt text "BTec_Logo-Orange"), whose "name"
#     attribute should read "image1.jpg" (it currently reads
#     "image2.jpg")
# This touches every header/footer in every section, for both the
# "first page" and "default" (primary) variants.

$d = $word.ActiveDocument

function Set-InlineShapeName($shape, $newName) {
    # A direct "$shape.Name = ..." assignment is unreliable on
    # picture shapes that live in a footer story in this host, so
    # route the write through the shape's own Selection (selecting
    # the shape first, then writing through $word.Selection) which
    # reliably commits the rename for both header and footer stories.
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($hf = 1; $hf -le 3; $hf++) {
        $header = $section.Headers.Item($hf)
        if ($header.Exists) {
            $count = $header.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shape = $header.Range.InlineShapes.Item($i)
                $alt = $shape.AlternativeText
                if ($alt -like "*PearsonLogo.png") {
                    Set-InlineShapeName $shape "image2.png"
                } elseif ($alt -eq "BTec_Logo-Orange") {
                    Set-InlineShapeName $shape "image1.jpg"
                }
            }
        }

        $footer = $section.Footers.Item($hf)
        if ($footer.Exists) {
            $count = $footer.Range.InlineShapes.Count
            for ($i = 1; $i -le $count; $i++) {
                $shape = $footer.Range.InlineShapes.Item($i)
                $alt = $shape.AlternativeText
                if ($alt -like "*PearsonLogo.png") {
                    Set-InlineShapeName $shape "image2.png"
                } elseif ($alt -eq "BTec_Logo-Orange") {
                    Set-InlineShapeName $shape "image1.jpg"
                }
            }
        }
    }
}

Write-Output "Renamed inline shapes."
